$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Treatment query SQL in cell B5 (TreatmentTab row) ---
# Bug fix: CONCAT(REPLACE(...)) simplified to REPLACE(...) for the
# "Treatment Agent" column.
$treatmentQuery = @"
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002599' AND srv.last_known_survival_status IN ('Alive')
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
"@

$b5 = $ws.Range("B5")
$b5.Value = $treatmentQuery

# Re-apply the existing wrap/size formatting; Excel allocates a fresh
# (but visually identical) style record for the cell during this kind
# of in-place edit.
$b5.WrapText = $true
$b5.Font.Size = 12
$b5.Font.ThemeColor = 1

# --- Update the saved view/selection state ---
$ws.Range("C5").Select()
